# Generate Report for Handback
# Adds the handback-freshness check results for d59b9702-3738-40e1-9d00-21648b2ac664
# to both the zh-cn and de-de worksheets (row 8 / column I,J,K,P), widens the
# "Error Detail" column, and links the new handback-file cell to its GitHub blob.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34d524ea5998888c2ee6c2d063f8fbfee0f617b6/e2e/d59b9702-3738-40e1-9d00-21648b2ac664.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb9d2cbfe37aff491360c7b5f3af0f55bd541c31/e2e/d59b9702-3738-40e1-9d00-21648b2ac664.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34d524ea5998888c2ee6c2d063f8fbfee0f617b6/e2e/d59b9702-3738-40e1-9d00-21648b2ac664.md."

# Hyperlink font color matching the workbook's existing "HyperLink" style (RGB 6495ED == BGR 15570276)
$hyperlinkColor = 15570276

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I8").Value = "d59b9702-3738-40e1-9d00-21648b2ac664.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestUrl, "", "", "d59b9702-3738-40e1-9d00-21648b2ac664.md")
$wsZh.Range("I8").Font.Color = $hyperlinkColor

$wsZh.Range("J8").Value = "d59b9702-3738-40e1-9d00-21648b2ac664.6f48dea37262a47703568b328ad140ec727c62e4.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-31 20:55:11"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I8").Value = "d59b9702-3738-40e1-9d00-21648b2ac664.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestUrl, "", "", "d59b9702-3738-40e1-9d00-21648b2ac664.md")
$wsDe.Range("I8").Font.Color = $hyperlinkColor

$wsDe.Range("J8").Value = "d59b9702-3738-40e1-9d00-21648b2ac664.6f48dea37262a47703568b328ad140ec727c62e4.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-31 20:55:23"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
